# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibition) — only "想去人数" (F column) counters refreshed.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 7692
$wsExpo.Range("F5").Value = 472
$wsExpo.Range("F6").Value = 4400
$wsExpo.Range("F8").Value = 620
$wsExpo.Range("F10").Value = 694
$wsExpo.Range("F11").Value = 172

# ---------------------------------------------------------------------------
# Sheet "演出" (Performance) — the 2024-08-03 concert ("菊次郎的夏天") is
# removed from the listing entirely; the CrossingX row shifts up to take its
# place (row 3 -> row 2) and its "想去人数" count ticks up from 14 to 15.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Rows.Item(2).Delete()
$wsShow.Range("A2").Value = 1
$wsShow.Range("F2").Value = 15

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) — same concert drops out of the combined feed;
# rows below shift up one position, their serial numbers (col A) need to be
# re-sequenced, and the refreshed "想去人数" counters applied.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Rows.Item(3).Delete()

for ($r = 3; $r -le 12; $r++) {
    $wsAll.Cells.Item($r, 1).Value = $r - 1
}

$wsAll.Range("F2").Value = 7692
$wsAll.Range("F5").Value = 472
$wsAll.Range("F6").Value = 4400
$wsAll.Range("F8").Value = 620
$wsAll.Range("F10").Value = 694
$wsAll.Range("F11").Value = 15
$wsAll.Range("F12").Value = 172

Write-Host "done"
